# Weekly fruit/vegetable price update: insert a new daily record for
# Espinaca (Vega Modelo de Temuco) before the existing row 16, shifting
# all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16; rows 16..82 shift down to 17..83.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row with the new day's data.
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44473
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 100112012
$ws.Cells.Item(16, 7).Value = "Espinaca"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 20
$ws.Cells.Item(16, 11).Value = 12000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 12000
$ws.Cells.Item(16, 14).Value = "$/docena de atados"
$ws.Cells.Item(16, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(16, 16).Value = 4000
$ws.Cells.Item(16, 17).Value = 3
$ws.Cells.Item(16, 18).Value = "Hortaliza"
